$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" '34.582.03'
Set-TextCell "E2" '  +1.64%  '

Set-TextCell "D3" '1.802.68'
Set-TextCell "E3" '  +1.08%  '

Set-TextCell "E4" '  -0.16%  '

Set-TextCell "D5" '224.10'
Set-TextCell "E5" '  -1.27%  '

Set-TextCell "E6" '  +0.08%  '

Set-TextCell "E7" '  -0.16%  '

Set-TextCell "D8" '32.42'
Set-TextCell "E8" '  +3.82%  '

Set-TextCell "E9" '  +3.08%  '

Set-TextCell "D10" '0.0711'
Set-TextCell "E10" '  +8.01%  '

Set-TextCell "E11" '  -0.04%  '

Set-TextCell "D12" '2.062.43'
Set-TextCell "E12" '  +1.05%  '

Set-TextCell "D13" '11.12'
Set-TextCell "E13" '  -2.21%  '

Set-TextCell "D14" '1.801.08'
Set-TextCell "E14" '  +0.89%  '

Set-TextCell "D15" '0.641'
Set-TextCell "E15" '  +1.35%  '

Set-TextCell "D16" '34.624.94'
Set-TextCell "E16" '  +1.69%  '

Set-TextCell "D17" '4.31'
Set-TextCell "E17" '  +2.17%  '

Set-TextCell "D18" '69.20'
Set-TextCell "E18" '  -0.19%  '

Set-TextCell "D19" '252.16'
Set-TextCell "E19" '  -0.15%  '

Set-TextCell "D20" '0.0₃0802'
Set-TextCell "E20" '  +8.23%  '

Set-TextCell "D21" '11.07'
Set-TextCell "E21" '  +6.24%  '

Set-TextCell "E22" '  -0.12%  '

Set-TextCell "E23" '  +0.26%  '

Set-TextCell "E24" '  +1.14%  '

Set-TextCell "D25" '161.66'
Set-TextCell "E25" '  +2.96%  '

Set-TextCell "D26" '16.40'
Set-TextCell "E26" '  -0.88%  '

Set-TextCell "D27" '7.13'
Set-TextCell "E27" '  +1.75%  '

Set-TextCell "E28" '  +0.10%  '

Set-TextCell "E29" '  -0.23%  '

Set-TextCell "D30" '580.82'
Set-TextCell "E30" '  +1,009.65%  '

Set-TextCell "E31" '  +2.09%  '

Set-TextCell "E32" '  +0.01%  '

Set-TextCell "E33" '  -0.28%  '

Set-TextCell "E34" '  +0.60%  '

Set-TextCell "E35" '  +2.57%  '

Set-TextCell "D36" '1.431.79'
Set-TextCell "E36" '  -1.24%  '

Set-TextCell "E37" '  +0.18%  '

Set-TextCell "D38" '0.642'
Set-TextCell "E38" '  +2.53%  '

Set-TextCell "E39" '  +2.74%  '

Set-TextCell "D40" '84.66'
Set-TextCell "E40" '  +1.62%  '

Set-TextCell "D41" '0.960'
Set-TextCell "E41" '  +7.00%  '

Set-TextCell "E42" '  -0.91%  '

Set-TextCell "E43" '  +0.15%  '

Set-TextCell "E44" '  +4.20%  '

Set-TextCell "D45" '6.05'
Set-TextCell "E45" '  +5.31%  '

Set-TextCell "E46" '  -0.85%  '

Set-TextCell "E47" '  -2.02%  '

Set-TextCell "B48" 'InjectiveProtocol'
Set-TextCell "C48" 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell "D48" '12.33'
Set-TextCell "E48" '  +3.82%  '

Set-TextCell "B49" 'RocketPoolETH'
Set-TextCell "C49" 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextCell "D49" '1.957.24'
Set-TextCell "E49" '  +0.81%  '

Set-TextCell "B50" 'Quant'
Set-TextCell "C50" 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell "D50" '106.72'
Set-TextCell "E50" '  +9.30%  '

Set-TextCell "E51" '  -0.06%  '
